$wb = $excel.ActiveWorkbook

# NOTE: "Vector_bf" and "Vector_BF" differ only by letter case, and
# Sheets.Item(name) resolves names case-insensitively (it would return the
# same sheet for both). Use the 1-based sheet index instead so each one is
# addressed unambiguously:
#   1 Funciones_Objetivo, 2 Restricciones_del_lider,
#   3 Restricciones_del_follower, 4 Punto_modificado, 5 Vector_bf,
#   6 Vector_BF, 7 Vector_Alpha
$wsFollower = $wb.Sheets.Item(3)
$wsPunto    = $wb.Sheets.Item(4)
$wsBf       = $wb.Sheets.Item(5)
$wsBF       = $wb.Sheets.Item(6)
$wsAlpha    = $wb.Sheets.Item(7)

# Cells that must hold digit-only-looking text need to be pre-formatted as
# Text ("@") before the value is assigned - otherwise Excel's COM layer
# auto-converts the numeric-looking string into a real number.  Apply the
# format once, up front, to exactly the cells that need it so every one of
# them shares the same single new style.
$textNumericRanges = @(
    $wsFollower.Range("B2:B6"),
    $wsFollower.Range("D2:F6"),
    $wsPunto.Range("A2:C2"),
    $wsBf.Range("A2:A3"),
    $wsBF.Range("A2:A4")
)
foreach ($rng in $textNumericRanges) {
    $rng.NumberFormat = "@"
}

# --- Restricciones_del_follower ---------------------------------------
# Row 2
$wsFollower.Range("A2").Value = "0.9261226324359138y_1 + 1.2878727323431254y_2"
$wsFollower.Range("B2").Value = "7.033987934013765"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.24011722556595838"
$wsFollower.Range("E2").Value = "0.07558256144476648"
$wsFollower.Range("F2").Value = "0.4651782448218301"

# Row 3
$wsFollower.Range("A3").Value = "-4 + 2.182352534978767y_1 + 0.790562119032984y_2"
$wsFollower.Range("B3").Value = "7.834642533158952"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.1083236165390392"
$wsFollower.Range("E3").Value = "0.2470951773333938"
$wsFollower.Range("F3").Value = "0.7055993983218097"

# Row 4 (Restriction_Set_Type unchanged: J_0_LP_v)
$wsFollower.Range("A4").Value = "-16 - 2x + 3.918163688312682y_1 - 0.7233534683711851y_2"
$wsFollower.Range("B4").Value = "-11.029067106621255"
$wsFollower.Range("D4").Value = "0.1102758390135593"
$wsFollower.Range("E4").Value = "0.7030203816296596"
$wsFollower.Range("F4").Value = "0.4629380383277565"

# Row 5
$wsFollower.Range("A5").Value = "-48 + 8x + 2.849178486278314y_1 + 1.2364251941225317y_2"
$wsFollower.Range("B5").Value = "14.885176000561037"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.3168885247170169"
$wsFollower.Range("E5").Value = "0.4962459620630362"
$wsFollower.Range("F5").Value = "0.9567530584293187"

# Row 6
$wsFollower.Range("A6").Value = "12 - 2x + 0.44656643466775314y_1 + 1.6358596000140058y_2"
$wsFollower.Range("B6").Value = "5.784706167175495"
$wsFollower.Range("C6").Value = "J_Ne_L0_v"
$wsFollower.Range("D6").Value = "0.4167665579899481"
$wsFollower.Range("E6").Value = "0.1244067936023273"
$wsFollower.Range("F6").Value = "0.40092091135914876"

# --- Punto_modificado ---------------------------------------------------
$wsPunto.Range("A2").Value = "5.875357499928848"
$wsPunto.Range("B2").Value = "4.657691821664619"
$wsPunto.Range("C2").Value = "2.112315956957238"

# --- Vector_bf ------------------------------------------------------------
$wsBf.Range("A2").Value = "-0.9798430273700574"
$wsBf.Range("A3").Value = "-1.3886890951177908"

# --- Vector_BF --------------------------------------------------------------
$wsBF.Range("A2").Value = "-1.3151133460403157"
$wsBF.Range("A3").Value = "-1.8332456561310178"
$wsBF.Range("A4").Value = "-2.6012356331310738"

# --- Vector_Alpha (stored as real numbers, not text) -----------------------
$wsAlpha.Range("A2").Value = 0.2602775476670889
$wsAlpha.Range("A3").Value = 0.17403064105926536
